# Update the "展览" and "全部类型" sheets with refreshed listing data.
# Both sheets share identical table contents in this workbook, and the
# commit regenerates the data for rows 2-5 and appends a new row 6.

$wb = $excel.ActiveWorkbook

$targetSheetNames = @("展览", "全部类型")

foreach ($sheetName in $targetSheetNames) {
    $ws = $null
    for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
        $candidate = $wb.Worksheets.Item($i)
        if ($candidate.Name -eq $sheetName) {
            $ws = $candidate
            break
        }
    }
    if ($ws -eq $null) {
        continue
    }

    # Make sure date-looking text (column B) and free-form text columns are
    # stored as literal text rather than being auto-converted to dates /
    # numbers by Excel's input parser.
    $textColumns = @(2, 3, 4, 5, 8, 9)
    foreach ($col in $textColumns) {
        $ws.Range($ws.Cells.Item(2, $col), $ws.Cells.Item(6, $col)).NumberFormat = "@"
    }

    # --- Row 2 ---
    $ws.Cells.Item(2, 2).Value = "2024-09-16"
    $ws.Cells.Item(2, 3).Value = "丽水·LZ栗子动漫游戏嘉年华（取消）"
    $ws.Cells.Item(2, 4).Value = "城北街798号 莱茵体育生活馆"
    $ws.Cells.Item(2, 5).Value = "2024.09.16 09:30-09.16 17:00"
    $ws.Cells.Item(2, 6).Value = 482
    $ws.Cells.Item(2, 7).NumberFormat = "@"
    $ws.Cells.Item(2, 7).Value = "不可售"
    $ws.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87480"
    $ws.Cells.Item(2, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/bATqcZhH1719285865931.jpeg"

    # --- Row 3 ---
    $ws.Cells.Item(3, 2).Value = "2024-10-01"
    $ws.Cells.Item(3, 3).Value = "丽水·CCAC动漫游戏嘉年华"
    $ws.Cells.Item(3, 4).Value = "南环西路109号 九城宴会中心"
    $ws.Cells.Item(3, 5).Value = "2024.10.01 09:00-10.01 16:00"
    $ws.Cells.Item(3, 6).Value = 101
    $ws.Cells.Item(3, 7).Value = 50
    $ws.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90985"
    $ws.Cells.Item(3, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/SasHidQZ1724379468667.jpeg"

    # --- Row 4 ---
    $ws.Cells.Item(4, 2).Value = "2024-10-01"
    $ws.Cells.Item(4, 3).Value = "丽水·熙梦动漫游戏展"
    $ws.Cells.Item(4, 4).Value = "城北街798号 莱茵体育生活馆"
    $ws.Cells.Item(4, 5).Value = "2024.10.01 10:00-10.01 17:00"
    $ws.Cells.Item(4, 6).Value = 4
    $ws.Cells.Item(4, 7).Value = 45
    $ws.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=92235"
    $ws.Cells.Item(4, 9).Value = "//i1.hdslb.com/bfs/openplatform/202409/JHGyuq6R1725938726636.jpeg"

    # --- Row 5 ---
    $ws.Cells.Item(5, 2).Value = "2024-10-02"
    $ws.Cells.Item(5, 3).Value = "青田·未闻展名国漫嘉年华"
    $ws.Cells.Item(5, 4).Value = "瓯南街道百悦城4幢 西娜君澜大饭店"
    $ws.Cells.Item(5, 5).Value = "2024.10.02 09:00-10.02 17:00"
    $ws.Cells.Item(5, 6).Value = 62
    $ws.Cells.Item(5, 7).Value = 45
    $ws.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91328"
    $ws.Cells.Item(5, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/w8uKtdlg1724147282076.jpeg"

    # --- Row 6 (new row) ---
    # Copy formatting from row 5's A cell (style index used for the row
    # number column) onto the new row's A cell before setting its value.
    $ws.Cells.Item(5, 1).Copy() | Out-Null
    $ws.Cells.Item(6, 1).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0

    $ws.Cells.Item(6, 1).Value = 5
    $ws.Cells.Item(6, 2).Value = "2024-10-03"
    $ws.Cells.Item(6, 3).Value = "缙云·星辰动漫游戏展嘉年华"
    $ws.Cells.Item(6, 4).Value = "黄龙路38号 中意大酒店(缙云店)"
    $ws.Cells.Item(6, 5).Value = "2024.10.03 10:00-10.03 17:00"
    $ws.Cells.Item(6, 6).Value = 11
    $ws.Cells.Item(6, 7).Value = 45
    $ws.Cells.Item(6, 8).Value = "https://show.bilibili.com/platform/detail.html?id=92236"
    $ws.Cells.Item(6, 9).Value = "//i0.hdslb.com/bfs/openplatform/202409/S13hVYA01725280725848.jpeg"
}

Write-Host "Updated sheets"
